$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-23 Saturday" "2024-03-24 Sunday"
Replace-Text "392×3=1176" "981×4=3924"
Replace-Text "568×7=3976" "809×8=6472"
Replace-Text "715×8=5720" "725×5=3625"
Replace-Text "367×8=2936" "598×3=1794"
Replace-Text "175×5=875" "986×4=3944"
Replace-Text "163×6=978" "275×7=1925"
Replace-Text "195×6=1170" "244×7=1708"
Replace-Text "558×9=5022" "948×8=7584"
Replace-Text "732×2=1464" "120×7=840"
Replace-Text "812×4=3248" "453×4=1812"
Replace-Text "106×3=318" "949×3=2847"
Replace-Text "262×4=1048" "883×4=3532"
Replace-Text "906×3=2718" "682×8=5456"
Replace-Text "545×2=1090" "137×3=411"
Replace-Text "519×6=3114" "639×5=3195"
Replace-Text "389×4=1556" "793×8=6344"
Replace-Text "190×4=760" "597×4=2388"
Replace-Text "619×3=1857" "828×5=4140"
Replace-Text "296×5=1480" "963×3=2889"
Replace-Text "491×2=982" "738×4=2952"
Replace-Text "672×4=2688" "548×4=2192"
Replace-Text "844×8=6752" "649×9=5841"
Replace-Text "429×5=2145" "527×4=2108"
Replace-Text "623×9=5607" "860×8=6880"
Replace-Text "537×6=3222" "433×9=3897"
